# Update "Avverkningsanmälningar" sheet:
#  1) Column C ("Förändrad") for every data row changes from 45184 to 45186.
#  2) For rows 2-12, the HYPERLINK formulas in columns S,T,V,W,X,Y gain a
#     second argument (the friendly display text), equal to the row's
#     "Beteckning" (column A) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$linkCols = @("S", "T", "V", "W", "X", "Y")

$lastRow = $ws.UsedRange.Rows.Count - 1
if ($lastRow -lt 262) {
    $lastRow = 262
}

for ($r = 2; $r -le $lastRow; $r++) {
    # 1) Update the "Förändrad" date value.
    $cCell = $ws.Range("C$r")
    if ($cCell.Value2 -eq 45184) {
        $cCell.Value = 45186
    }

    # 2) Rewrite hyperlink formulas (only present on rows 2-12) so the
    #    HYPERLINK() call includes the friendly-name second argument.
    $beteckning = $ws.Range("A$r").Value2

    if ($beteckning) {
        foreach ($col in $linkCols) {
            $cell = $ws.Range("$col$r")
            $formula = $cell.Formula
            if ($formula -and $formula.Length -gt 0) {
                $trimmed = $formula.TrimEnd()
                if ($trimmed.EndsWith(")") -and -not ($trimmed.Contains(', "'))) {
                    $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $beteckning + '")'
                    $cell.Formula = $newFormula
                }
            }
        }
    }
}
